# Update LR-pair TPM values for Tnf-Tnfrsf21 (rows 2-26, columns A-T).
# Values reflect the refreshed NATMI TPM run described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "ECs"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 5.429001
$arr[0,7] = 16.287003
$arr[0,8] = 0.04230716253661782
$arr[0,9] = 0.04239440107683373
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 4.631270333333333
$arr[0,13] = 13.893811
$arr[0,14] = 0.05846361049715151
$arr[0,15] = 0.0596002562356855
$arr[0,16] = 25.143171270937
$arr[0,17] = 226.288541438433
$arr[0,18] = 0.002473429471780504
$arr[0,19] = 0.002526717167137711
$ws.Range("A2:T2").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "ECs"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "FAPs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 5.429001
$arr[0,7] = 16.287003
$arr[0,8] = 0.04230716253661782
$arr[0,9] = 0.04239440107683373
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 1.399706666666667
$arr[0,13] = 4.199120000000001
$arr[0,14] = 0.01766942965546306
$arr[0,15] = 0.01801295756537869
$arr[0,16] = 7.599008893040001
$arr[0,17] = 68.39108003736
$arr[0,18] = 0.0007475434323630109
$arr[0,19] = 0.0007636485476066507
$ws.Range("A3:T3").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "ECs"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "Inflammatory-Mac"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 5.429001
$arr[0,7] = 16.287003
$arr[0,8] = 0.04230716253661782
$arr[0,9] = 0.04239440107683373
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 33.32967466666667
$arr[0,13] = 99.989024
$arr[0,14] = 0.4207426855832669
$arr[0,15] = 0.428922737696382
$arr[0,16] = 180.946837095008
$arr[0,17] = 1628.521533855072
$arr[0,18] = 0.01780042918506436
$arr[0,19] = 0.01818392257287397
$ws.Range("A4:T4").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "ECs"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "MuSCs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 5.429001
$arr[0,7] = 16.287003
$arr[0,8] = 0.04230716253661782
$arr[0,9] = 0.04239440107683373
$arr[0,10] = 2
$arr[0,11] = 1
$arr[0,12] = 4.5322385
$arr[0,13] = 9.064477
$arr[0,14] = 0.05721346569581108
$arr[0,15] = 0.03888387079991788
$arr[0,16] = 24.6055273487385
$arr[0,17] = 147.633164092431
$arr[0,18] = 0.002420539392475887
$arr[0,19] = 0.001648458414111502
$ws.Range("A5:T5").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "ECs"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "Resolving-Mac"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 5.429001
$arr[0,7] = 16.287003
$arr[0,8] = 0.04230716253661782
$arr[0,9] = 0.04239440107683373
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 35.3234
$arr[0,13] = 105.9702
$arr[0,14] = 0.4459108085683075
$arr[0,15] = 0.454580177702636
$arr[0,16] = 191.7707739234
$arr[0,17] = 1725.9369653106
$arr[0,18] = 0.01886522105493406
$arr[0,19] = 0.0192716543751039
$ws.Range("A6:T6").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "FAPs"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 1.36117
$arr[0,7] = 4.08351
$arr[0,8] = 0.01060733649339319
$arr[0,9] = 0.01062920911485442
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 4.631270333333333
$arr[0,13] = 13.893811
$arr[0,14] = 0.05846361049715151
$arr[0,15] = 0.0596002562356855
$arr[0,16] = 6.303946239623334
$arr[0,17] = 56.73551615661
$arr[0,18] = 0.0006201431891619602
$arr[0,19] = 0.0006335035868280074
$ws.Range("A7:T7").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "FAPs"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "FAPs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 1.36117
$arr[0,7] = 4.08351
$arr[0,8] = 0.01060733649339319
$arr[0,9] = 0.01062920911485442
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 1.399706666666667
$arr[0,13] = 4.199120000000001
$arr[0,14] = 0.01766942965546306
$arr[0,15] = 0.01801295756537869
$arr[0,16] = 1.905238723466667
$arr[0,17] = 17.1471485112
$arr[0,18] = 0.0001874255860018371
$arr[0,19] = 0.0001914634927394091
$ws.Range("A8:T8").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "FAPs"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "Inflammatory-Mac"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 1.36117
$arr[0,7] = 4.08351
$arr[0,8] = 0.01060733649339319
$arr[0,9] = 0.01062920911485442
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 33.32967466666667
$arr[0,13] = 99.989024
$arr[0,14] = 0.4207426855832669
$arr[0,15] = 0.428922737696382
$arr[0,16] = 45.36735326602668
$arr[0,17] = 408.30617939424
$arr[0,18] = 0.004462959243115642
$arr[0,19] = 0.004559109473090696
$ws.Range("A9:T9").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "FAPs"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "MuSCs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 1.36117
$arr[0,7] = 4.08351
$arr[0,8] = 0.01060733649339319
$arr[0,9] = 0.01062920911485442
$arr[0,10] = 2
$arr[0,11] = 1
$arr[0,12] = 4.5322385
$arr[0,13] = 9.064477
$arr[0,14] = 0.05721346569581108
$arr[0,15] = 0.03888387079991788
$arr[0,16] = 6.169147079045001
$arr[0,17] = 37.01488247427
$arr[0,18] = 0.000606882482588676
$arr[0,19] = 0.0004133047939273088
$ws.Range("A10:T10").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "FAPs"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "Resolving-Mac"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 1.36117
$arr[0,7] = 4.08351
$arr[0,8] = 0.01060733649339319
$arr[0,9] = 0.01062920911485442
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 35.3234
$arr[0,13] = 105.9702
$arr[0,14] = 0.4459108085683075
$arr[0,15] = 0.454580177702636
$arr[0,16] = 48.08115237800001
$arr[0,17] = 432.7303714020001
$arr[0,18] = 0.004729925992525071
$arr[0,19] = 0.004831827768269002
$ws.Range("A11:T11").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "Inflammatory-Mac"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 63.64001366666667
$arr[0,7] = 190.920041
$arr[0,8] = 0.4959344089323702
$arr[0,9] = 0.4969570394110899
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 4.631270333333333
$arr[0,13] = 13.893811
$arr[0,14] = 0.05846361049715151
$arr[0,15] = 0.0596002562356855
$arr[0,16] = 294.7341073073612
$arr[0,17] = 2652.606965766251
$arr[0,18] = 0.02899411611595714
$arr[0,19] = 0.02961876688702861
$ws.Range("A12:T12").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "Inflammatory-Mac"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "FAPs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 63.64001366666667
$arr[0,7] = 190.920041
$arr[0,8] = 0.4959344089323702
$arr[0,9] = 0.4969570394110899
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 1.399706666666667
$arr[0,13] = 4.199120000000001
$arr[0,14] = 0.01766942965546306
$arr[0,15] = 0.01801295756537869
$arr[0,16] = 89.07735139599113
$arr[0,17] = 801.6961625639201
$arr[0,18] = 0.008762878152354167
$arr[0,19] = 0.008951666062728191
$ws.Range("A13:T13").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "Inflammatory-Mac"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "Inflammatory-Mac"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 63.64001366666667
$arr[0,7] = 190.920041
$arr[0,8] = 0.4959344089323702
$arr[0,9] = 0.4969570394110899
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 33.32967466666667
$arr[0,13] = 99.989024
$arr[0,14] = 0.4207426855832669
$arr[0,15] = 0.428922737696382
$arr[0,16] = 2121.100951292221
$arr[0,17] = 19089.90856162999
$arr[0,18] = 0.2086607750873555
$arr[0,19] = 0.2131561738616935
$ws.Range("A14:T14").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "Inflammatory-Mac"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "MuSCs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 63.64001366666667
$arr[0,7] = 190.920041
$arr[0,8] = 0.4959344089323702
$arr[0,9] = 0.4969570394110899
$arr[0,10] = 2
$arr[0,11] = 1
$arr[0,12] = 4.5322385
$arr[0,13] = 9.064477
$arr[0,14] = 0.05721346569581108
$arr[0,15] = 0.03888387079991788
$arr[0,16] = 288.4317200805929
$arr[0,17] = 1730.590320483557
$arr[0,18] = 0.0283741262928245
$arr[0,19] = 0.01932361331357052
$ws.Range("A15:T15").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "Inflammatory-Mac"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "Resolving-Mac"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 63.64001366666667
$arr[0,7] = 190.920041
$arr[0,8] = 0.4959344089323702
$arr[0,9] = 0.4969570394110899
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 35.3234
$arr[0,13] = 105.9702
$arr[0,14] = 0.4459108085683075
$arr[0,15] = 0.454580177702636
$arr[0,16] = 2247.981658753133
$arr[0,17] = 20231.8349287782
$arr[0,18] = 0.2211425132838788
$arr[0,19] = 0.2259068192860691
$ws.Range("A16:T16").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "MuSCs"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "ECs"
$arr[0,4] = 1
$arr[0,5] = 0.5
$arr[0,6] = 0.792186
$arr[0,7] = 1.584372
$arr[0,8] = 0.006173353414603005
$arr[0,9] = 0.004124055360148531
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 4.631270333333333
$arr[0,13] = 13.893811
$arr[0,14] = 0.05846361049715151
$arr[0,15] = 0.0596002562356855
$arr[0,16] = 3.668827520282
$arr[0,17] = 22.012965121692
$arr[0,18] = 0.0003609165294926103
$arr[0,19] = 0.0002457947561950046
$ws.Range("A17:T17").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "MuSCs"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "FAPs"
$arr[0,4] = 1
$arr[0,5] = 0.5
$arr[0,6] = 0.792186
$arr[0,7] = 1.584372
$arr[0,8] = 0.006173353414603005
$arr[0,9] = 0.004124055360148531
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 1.399706666666667
$arr[0,13] = 4.199120000000001
$arr[0,14] = 0.01766942965546306
$arr[0,15] = 0.01801295756537869
$arr[0,16] = 1.10882802544
$arr[0,17] = 6.652968152640001
$arr[0,18] = 0.0001090796338976405
$arr[0,19] = 0.00007428643419962803
$ws.Range("A18:T18").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "MuSCs"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "Inflammatory-Mac"
$arr[0,4] = 1
$arr[0,5] = 0.5
$arr[0,6] = 0.792186
$arr[0,7] = 1.584372
$arr[0,8] = 0.006173353414603005
$arr[0,9] = 0.004124055360148531
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 33.32967466666667
$arr[0,13] = 99.989024
$arr[0,14] = 0.4207426855832669
$arr[0,15] = 0.428922737696382
$arr[0,16] = 26.403301655488
$arr[0,17] = 158.419809932928
$arr[0,18] = 0.002597393294714699
$arr[0,19] = 0.001768901115486346
$ws.Range("A19:T19").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "MuSCs"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "MuSCs"
$arr[0,4] = 1
$arr[0,5] = 0.5
$arr[0,6] = 0.792186
$arr[0,7] = 1.584372
$arr[0,8] = 0.006173353414603005
$arr[0,9] = 0.004124055360148531
$arr[0,10] = 2
$arr[0,11] = 1
$arr[0,12] = 4.5322385
$arr[0,13] = 9.064477
$arr[0,14] = 0.05721346569581108
$arr[0,15] = 0.03888387079991788
$arr[0,16] = 3.590375888361
$arr[0,17] = 14.361503553444
$arr[0,18] = 0.0003531989438145072
$arr[0,19] = 0.0001603592357957243
$ws.Range("A20:T20").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "MuSCs"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "Resolving-Mac"
$arr[0,4] = 1
$arr[0,5] = 0.5
$arr[0,6] = 0.792186
$arr[0,7] = 1.584372
$arr[0,8] = 0.006173353414603005
$arr[0,9] = 0.004124055360148531
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 35.3234
$arr[0,13] = 105.9702
$arr[0,14] = 0.4459108085683075
$arr[0,15] = 0.454580177702636
$arr[0,16] = 27.9827029524
$arr[0,17] = 167.8962177144
$arr[0,18] = 0.002752765012683548
$arr[0,19] = 0.001874713818471828
$ws.Range("A21:T21").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "Resolving-Mac"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 57.101078
$arr[0,7] = 171.303234
$arr[0,8] = 0.444977738623016
$arr[0,9] = 0.4458952950370734
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 4.631270333333333
$arr[0,13] = 13.893811
$arr[0,14] = 0.05846361049715151
$arr[0,15] = 0.0596002562356855
$arr[0,16] = 264.4505285427526
$arr[0,17] = 2380.054756884774
$arr[0,18] = 0.0260150051907593
$arr[0,19] = 0.02657547383849616
$ws.Range("A22:T22").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "Resolving-Mac"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "FAPs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 57.101078
$arr[0,7] = 171.303234
$arr[0,8] = 0.444977738623016
$arr[0,9] = 0.4458952950370734
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 1.399706666666667
$arr[0,13] = 4.199120000000001
$arr[0,14] = 0.01766942965546306
$arr[0,15] = 0.01801295756537869
$arr[0,16] = 79.92475955045335
$arr[0,17] = 719.3228359540801
$arr[0,18] = 0.00786250285084641
$arr[0,19] = 0.008031893028104817
$ws.Range("A23:T23").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "Resolving-Mac"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "Inflammatory-Mac"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 57.101078
$arr[0,7] = 171.303234
$arr[0,8] = 0.444977738623016
$arr[0,9] = 0.4458952950370734
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 33.32967466666667
$arr[0,13] = 99.989024
$arr[0,14] = 0.4207426855832669
$arr[0,15] = 0.428922737696382
$arr[0,16] = 1903.160352855958
$arr[0,17] = 17128.44317570362
$arr[0,18] = 0.1872211287730167
$arr[0,19] = 0.1912546306732375
$ws.Range("A24:T24").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "Resolving-Mac"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "MuSCs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 57.101078
$arr[0,7] = 171.303234
$arr[0,8] = 0.444977738623016
$arr[0,9] = 0.4458952950370734
$arr[0,10] = 2
$arr[0,11] = 1
$arr[0,12] = 4.5322385
$arr[0,13] = 9.064477
$arr[0,14] = 0.05721346569581108
$arr[0,15] = 0.03888387079991788
$arr[0,16] = 258.795704103103
$arr[0,17] = 1552.774224618618
$arr[0,18] = 0.02545871858410751
$arr[0,19] = 0.01733813504251282
$ws.Range("A25:T25").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "Resolving-Mac"
$arr[0,1] = "Tnf"
$arr[0,2] = "Tnfrsf21"
$arr[0,3] = "Resolving-Mac"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 57.101078
$arr[0,7] = 171.303234
$arr[0,8] = 0.444977738623016
$arr[0,9] = 0.4458952950370734
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 35.3234
$arr[0,13] = 105.9702
$arr[0,14] = 0.4459108085683075
$arr[0,15] = 0.454580177702636
$arr[0,16] = 2017.0042186252
$arr[0,17] = 18153.0379676268
$arr[0,18] = 0.198420383224286
$arr[0,19] = 0.2026951624547222
$ws.Range("A26:T26").Value = $arr
